$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 17 (allcount): extend "Used by" list to include SelectGermlineVariants, RegionalExpression
$ws.Range("I17").Value = "ExpressionDistribution, SelectGermlineVariants, RegionalExpression"

# 2) Row 34 (<analysis>.selectedVariants): extend "Using Inputs" to include allcount (DNA & RNA)
$ws.Range("G34").Value = "*.vcf, experiments.txt, allcount (DNA & RNA)"

# 3) Insert a new row at 36 for <analysis>.annotatedSelectedVariants.txt
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "<analysis>.annotatedSelectedVariants.txt"
$ws.Range("B36").Value = "text"
$ws.Range("C36").Value = "M"
$ws.Range("D36").Value = "selected variants + RNA & DNA read counts at the variant locus"
$ws.Range("E36").Value = "mutant-expression"
$ws.Range("F36").Value = "\tcga\<disease>\tumor\normal\{wgs|wxs}\<analysis-id>\<analysis-id>.annotatedSelectedVariants"
$ws.Range("G36").Value = "indices, *.selectedVariants, *.readsAtSelectedVariants"

# 4) Fill in new rows 39-42 for the ExpressionByMutationCount pipeline (previously blank filler rows)
# row 39
$ws.Range("A39").Value = "<gene>_lines.txt"
$ws.Range("B39").Value = "text"
$ws.Range("C39").Value = "M"
$ws.Range("D39").Value = "Normalized expression level around genes & mutation count for each tumor"
$ws.Range("E39").Value = "ExpressionByMutationCount"
$ws.Range("F39").Value = "\temp\expression\RegionalExpressionByGene"
$ws.Range("G39").Value = "experiments.txt, *.gene_expression.txt, <gene>.unfiltered_counts.txt"
$ws.Range("I39").Value = "Final results"

# row 40
$ws.Range("A40").Value = "ExpressionDistributionByMutationCount"
$ws.Range("B40").Value = "text"
$ws.Range("C40").Value = "K"
$ws.Range("D40").Value = "p values for expression around genes by mutation count"
$ws.Range("E40").Value = "ExpressionByMutationCount"
$ws.Range("F40").Value = "\temp\expression\"
$ws.Range("G40").Value = "experiments.txt, *.gene_expression.txt, <gene>.unfiltered_counts.txt"
$ws.Range("I40").Value = "Final results"

# row 41
$ws.Range("A41").Value = "<gene>_allele_specific_lines.txt"
$ws.Range("B41").Value = "text"
$ws.Range("C41").Value = "M"
$ws.Range("D41").Value = "Normalized allele-specific expression level around genes & mutation count for each tumor"
$ws.Range("E41").Value = "ExpressionByMutationCount"
$ws.Range("F41").Value = "\temp\expression\RegionalExpressionByGene"
$ws.Range("G41").Value = "experiments.txt, *.allele_specific_gene_expression.txt, <gene>.unfiltered_counts.txt"
$ws.Range("I41").Value = "Final results"

# row 42
$ws.Range("A42").Value = "AlleleSpecificExpressionDistributionByMutationCount"
$ws.Range("B42").Value = "text"
$ws.Range("C42").Value = "K"
$ws.Range("D42").Value = "p values for allele-specific expression around genes by mutation count"
$ws.Range("E42").Value = "ExpressionByMutationCount"
$ws.Range("F42").Value = "\temp\expression\"
$ws.Range("G42").Value = "experiments.txt, *.allele_specific_gene_expression.txt, <gene>.unfiltered_counts.txt"
$ws.Range("I42").Value = "Final results"

# 5) Column widths (closest achievable via ColumnWidth quantization)
$ws.Columns.Item(1).ColumnWidth = 49.8
$ws.Columns.Item(4).ColumnWidth = 80.3
$ws.Columns.Item(6).ColumnWidth = 87.3
$ws.Columns.Item(7).ColumnWidth = 76.3

# 6) Update selection to match author (G47) and re-select the active sheet
$ws.Range("G47").Select()
